$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 459.2857
$ws.Range("I4").Value = 459.2857
$ws.Range("K4").Value = 459.2857
$ws.Range("M4").Value = -345.2857
$ws.Range("H6").Value = 71428760
$ws.Range("I6").Value = 100000130
$ws.Range("K6").Value = 300000390
$ws.Range("M6").Value = -300000278
$ws.Range("H38").Value = 9908.223
$ws.Range("I38").Value = 22191.8
$ws.Range("J38").Value = 5183.769
$ws.Range("K38").Value = 66575.39999999999
$ws.Range("L38").Value = 15551.307
$ws.Range("M38").Value = -66203.39999999999
$ws.Range("N38").Value = -16295.307
$ws.Range("H41").Value = 558.3333
$ws.Range("I41").Value = 649.1667
$ws.Range("J41").Value = 376.66666
$ws.Range("K41").Value = 649.1667
$ws.Range("L41").Value = 376.66666
$ws.Range("M41").Value = -209.1667
$ws.Range("N41").Value = -1256.66666
$ws.Range("H51").Value = 4399.8
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H57").Value = 64582.668
$ws.Range("J57").Value = 64582.668
$ws.Range("L57").Value = 193748.004
$ws.Range("N57").Value = -194746.004
$ws.Range("H100").Value = 3293.6667
$ws.Range("I100").Value = 1963.9231
$ws.Range("K100").Value = 1963.9231
$ws.Range("M100").Value = -1422.9231
$ws.Range("H112").Value = 1383
$ws.Range("I112").Value = 1200
$ws.Range("K112").Value = 3600
$ws.Range("M112").Value = -2492
$ws.Range("H129").Value = 1971.375
$ws.Range("J129").Value = 2554.4
$ws.Range("L129").Value = 7663.200000000001
$ws.Range("N129").Value = -17663.2
$ws.Range("H133").Value = 104996
$ws.Range("J133").Value = 104996
$ws.Range("L133").Value = 104996
$ws.Range("N133").Value = -115116
$ws.Range("H136").Value = 137390.72
$ws.Range("J136").Value = 137390.72
$ws.Range("L136").Value = 137390.72
$ws.Range("N136").Value = -147590.72
$ws.Range("H138").Value = 3591674
$ws.Range("I138").Value = 14848.5
$ws.Range("K138").Value = 44545.5
$ws.Range("M138").Value = -39405.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4352.058
$ws.Range("I32").Value = 3337.5454
$ws.Range("K32").Value = 3337.5454
$ws.Range("M32").Value = -3050.5454
$ws.Range("H45").Value = 17517
$ws.Range("I45").Value = 34873
$ws.Range("K45").Value = 34873
$ws.Range("M45").Value = -34496
$ws.Range("H88").Value = 2612.8
$ws.Range("I88").Value = 2743.2856
$ws.Range("J88").Value = 2498.625
$ws.Range("K88").Value = 2743.2856
$ws.Range("L88").Value = 2498.625
$ws.Range("M88").Value = -2337.2856
$ws.Range("N88").Value = -3310.625
$ws.Range("H91").Value = 2612.8
$ws.Range("I91").Value = 2743.2856
$ws.Range("J91").Value = 2498.625
$ws.Range("K91").Value = 2743.2856
$ws.Range("L91").Value = 2498.625
$ws.Range("M91").Value = -1339.2856
$ws.Range("N91").Value = -5306.625
$ws.Range("H102").Value = 4735
$ws.Range("I102").Value = 3682
$ws.Range("K102").Value = 3682
$ws.Range("M102").Value = -2060
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6894.9
$ws.Range("I105").Value = 8354.143
$ws.Range("K105").Value = 8354.143
$ws.Range("M105").Value = -6607.143
$ws.Range("H134").Value = 2299.279
$ws.Range("I134").Value = 2131.5945
$ws.Range("K134").Value = 6394.7835
$ws.Range("M134").Value = -3859.7835
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 18290.75
$ws.Range("J28").Value = 18290.75
$ws.Range("L28").Value = 18290.75
$ws.Range("N28").Value = -18780.75
$ws.Range("H47").Value = 13499.5
$ws.Range("I47").Value = 11999
$ws.Range("K47").Value = 11999
$ws.Range("M47").Value = -11433
$ws.Range("H114").Value = 16999.844
$ws.Range("J114").Value = 19358.975
$ws.Range("L114").Value = 19358.975
$ws.Range("N114").Value = -28036.975
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7717.5454
$ws.Range("I56").Value = 7717.5454
$ws.Range("K56").Value = 7717.5454
$ws.Range("M56").Value = -7187.5454
$ws.Range("H94").Value = 12948.75
$ws.Range("I94").Value = 2897.5
$ws.Range("J94").Value = 23000
$ws.Range("K94").Value = 8692.5
$ws.Range("L94").Value = 69000
$ws.Range("M94").Value = -8016.5
$ws.Range("N94").Value = -70352
$ws.Range("H97").Value = 1534.5714
$ws.Range("J97").Value = 2749.5
$ws.Range("L97").Value = 8248.5
$ws.Range("N97").Value = -9240.5
$ws.Range("H107").Value = 2888
$ws.Range("I107").Value = 3396.3333
$ws.Range("J107").Value = 2379.6667
$ws.Range("K107").Value = 10188.9999
$ws.Range("L107").Value = 7139.000100000001
$ws.Range("M107").Value = -8268.999899999999
$ws.Range("N107").Value = -10979.0001
$ws.Range("H128").Value = 285015
$ws.Range("I128").Value = 285015
$ws.Range("K128").Value = 855045
$ws.Range("M128").Value = -850065
$ws.Range("H131").Value = 2670365.5
$ws.Range("I131").Value = 1210.125
$ws.Range("J131").Value = 3926438.8
$ws.Range("K131").Value = 3630.375
$ws.Range("L131").Value = 11779316.4
$ws.Range("M131").Value = 1409.625
$ws.Range("N131").Value = -11789396.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19627776
$ws.Range("I70").Value = 30321748
$ws.Range("J70").Value = 22164.834
$ws.Range("K70").Value = 30321748
$ws.Range("L70").Value = 22164.834
$ws.Range("M70").Value = -30321478
$ws.Range("N70").Value = -22704.834
$ws.Range("H73").Value = 19627776
$ws.Range("I73").Value = 30321748
$ws.Range("J73").Value = 22164.834
$ws.Range("K73").Value = 30321748
$ws.Range("L73").Value = 22164.834
$ws.Range("M73").Value = -30320812
$ws.Range("N73").Value = -24036.834
$ws.Range("H126").Value = 2913.5
$ws.Range("J126").Value = 3854
$ws.Range("L126").Value = 11562
$ws.Range("N126").Value = -16502
$ws.Range("H129").Value = 50000
$ws.Range("J129").Value = 50000
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9959.615
$ws.Range("I16").Value = 889.5454999999999
$ws.Range("J16").Value = 59845
$ws.Range("K16").Value = 889.5454999999999
$ws.Range("L16").Value = 59845
$ws.Range("M16").Value = -719.5454999999999
$ws.Range("N16").Value = -60185
$ws.Range("H22").Value = 3514.7896
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 3654.5
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 3654.5
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -4244.5
$ws.Range("H27").Value = 3514.7896
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 3654.5
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 3654.5
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -3868.5
$ws.Range("H55").Value = 1617.625
$ws.Range("J55").Value = 2160.3
$ws.Range("L55").Value = 2160.3
$ws.Range("N55").Value = -2506.3
$ws.Range("H68").Value = 6178.8335
$ws.Range("I68").Value = 3527.7144
$ws.Range("J68").Value = 7270.4707
$ws.Range("K68").Value = 3527.7144
$ws.Range("L68").Value = 7270.4707
$ws.Range("M68").Value = -2778.7144
$ws.Range("N68").Value = -8768.4707
$ws.Range("H71").Value = 6178.8335
$ws.Range("I71").Value = 3527.7144
$ws.Range("J71").Value = 7270.4707
$ws.Range("K71").Value = 17638.572
$ws.Range("L71").Value = 36352.3535
$ws.Range("M71").Value = -13894.572
$ws.Range("N71").Value = -43840.3535
$ws.Range("H93").Value = 4243.2188
$ws.Range("I93").Value = 1848.7333
$ws.Range("K93").Value = 1848.7333
$ws.Range("M93").Value = -600.7333000000001
$ws.Range("H120").Value = 202000
$ws.Range("J120").Value = 202000
$ws.Range("L120").Value = 202000
$ws.Range("N120").Value = -211676
$ws.Range("H122").Value = 3544.7932
$ws.Range("I122").Value = 3381.4443
$ws.Range("K122").Value = 10144.3329
$ws.Range("M122").Value = -7694.332900000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 5666.3335
$ws.Range("J4").Value = 5666.3335
$ws.Range("L4").Value = 5666.3335
$ws.Range("N4").Value = -5892.3335
$ws.Range("H96").Value = 3671.5334
$ws.Range("I96").Value = 2226.4285
$ws.Range("J96").Value = 4936
$ws.Range("K96").Value = 2226.4285
$ws.Range("L96").Value = 4936
$ws.Range("M96").Value = -853.4285
$ws.Range("N96").Value = -7682
$ws.Range("H132").Value = 8205.162
$ws.Range("I132").Value = 8737.971
$ws.Range("J132").Value = 2166.6667
$ws.Range("K132").Value = 26213.913
$ws.Range("L132").Value = 6500.000100000001
$ws.Range("M132").Value = -23683.913
$ws.Range("N132").Value = -11560.0001
$ws.Range("H136").Value = 3889.1064
$ws.Range("I136").Value = 2733.5
$ws.Range("J136").Value = 8768.333000000001
$ws.Range("K136").Value = 8200.5
$ws.Range("L136").Value = 26304.999
$ws.Range("M136").Value = -5650.5
$ws.Range("N136").Value = -31404.999
